# Updated cryptos list with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto table on Sheet1, row by row, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price text (or $null if unchanged),
# new Volume(1h) text (or $null if unchanged).
$updates = @(
    @{ Row = 2;  D = "46.197.43";  E = "  +3.83%  " },
    @{ Row = 3;  D = "2.454.16";   E = "  +1.40%  " },
    @{ Row = 4;  D = $null;        E = "  +0.01%  " },
    @{ Row = 5;  D = "320.77";     E = "  +2.07%  " },
    @{ Row = 6;  D = "105.39";     E = "  +4.75%  " },
    @{ Row = 7;  D = "0.517";      E = "  +0.93%  " },
    @{ Row = 8;  D = $null;        E = "  +0.00%  " },
    @{ Row = 9;  D = "0.535";      E = "  +2.45%  " },
    @{ Row = 10; D = "36.11";      E = "  +2.41%  " },
    @{ Row = 11; D = "0.0818";     E = "  +2.45%  " },
    @{ Row = 12; D = "0.122";      E = "  +0.57%  " },
    @{ Row = 13; D = "18.39";      E = "  -4.23%  " },
    @{ Row = 14; D = "7.10";       E = $null },
    @{ Row = 15; D = "2.846.71";   E = "  +1.68%  " },
    @{ Row = 16; D = "2.474.33";   E = "  +1.85%  " },
    @{ Row = 17; D = "0.842";      E = "  +1.40%  " },
    @{ Row = 18; D = "46.096.78";  E = "  +3.90%  " },
    @{ Row = 19; D = $null;        E = "  +2.28%  " },
    @{ Row = 20; D = "6.43";       E = "  +0.65%  " },
    @{ Row = 21; D = "0.0₃0939";   E = "  +2.54%  " },
    @{ Row = 22; D = "70.93";      E = "  +3.22%  " },
    @{ Row = 23; D = $null;        E = "  +4.58%  " },
    @{ Row = 24; D = "247.63";     E = "  +2.16%  " },
    @{ Row = 25; D = $null;        E = "  +1.86%  " },
    @{ Row = 26; D = "25.96";      E = "  +3.35%  " },
    @{ Row = 27; D = $null;        E = "  -0.02%  " },
    @{ Row = 28; D = "2.29";       E = "  +0.19%  " },
    @{ Row = 29; D = "9.73";       E = "  +1.62%  " },
    @{ Row = 30; D = "34.74";      E = "  +4.80%  " },
    @{ Row = 31; D = "49.19";      E = "  +1.65%  " },
    @{ Row = 32; D = $null;        E = "  +4.85%  " },
    @{ Row = 33; D = "19.72";      E = "  +2.84%  " },
    @{ Row = 34; D = $null;        E = "  +3.48%  " },
    @{ Row = 35; D = $null;        E = "  +0.00%  " },
    @{ Row = 36; D = "0.0766";     E = "  -0.85%  " },
    @{ Row = 37; D = "4.57";       E = "  +2.26%  " },
    @{ Row = 39; D = $null;        E = "  +2.91%  " },
    @{ Row = 40; D = "125.45";     E = "  +4.72%  " },
    @{ Row = 41; D = $null;        E = "  +1.82%  " },
    @{ Row = 42; D = $null;        E = "  +0.58%  " },
    @{ Row = 43; D = "20.78";      E = "  -0.37%  " },
    @{ Row = 44; D = "0.0292";     E = "  +1.29%  " },
    @{ Row = 45; D = "1.972.68";   E = "  +1.54%  " },
    @{ Row = 46; D = $null;        E = "  +1.51%  " },
    @{ Row = 47; D = $null;        E = "  -4.02%  " },
    @{ Row = 48; D = $null;        E = "  +12.27%  " },
    @{ Row = 50; D = "5.07";       E = "  +9.73%  " },
    @{ Row = 51; D = "78.29";      E = "  +5.07%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D = Price

        # Some prices (e.g. "320.77") look like plain numbers, and Excel
        # would silently convert them to a floating point value on
        # assignment, losing the exact original text (and trailing zero
        # formatting, e.g. "7.10" -> 7.1). Force the cell to Text format
        # first in that case so the literal string is preserved, matching
        # the original inline-string cell contents. Values that already
        # contain multiple dots or other non-numeric characters (e.g.
        # "2.846.71", "0.0₃0939") are naturally stored as text, so leave
        # their formatting untouched.
        $looksNumeric = $u.D -match '^[+-]?\d+(\.\d+)?$'
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E   # column E = Volume(1h)
    }
}
